$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Summary" sheet - refresh aggregate metrics
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1198.85   # Current Capital
$summary.Range("B4").Value = -1.15    # Total P&L $
$summary.Range("B5").Value = -0.7     # Total P&L %
$summary.Range("B6").Value = 33       # Total Trades
$summary.Range("B7").Value = 13       # Winning Trades
$summary.Range("B9").Value = 39.39    # Win Rate %

# ---------------------------------------------------------------------------
# 2) "Strategy Status" sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 98.84999999999999  # Capital
$status.Range("D4").Value = 33                 # Trades
$status.Range("E4").Value = -1.15              # P&L $
$status.Range("F4").Value = -1.15              # P&L %
$status.Range("G4").Value = 39.39              # Win Rate %

# ---------------------------------------------------------------------------
# 3) Append the newly closed trade (#33) to both the "All Trades" and the
#    "MarketMaking" sheets as row 34.
# ---------------------------------------------------------------------------
function Add-TradeRow34($ws) {
    # Text-like columns (Date / Time) must stay as plain text, not get
    # auto-converted into date/time serial values. Temporarily force a text
    # number format, write the value, then clear the format again so the
    # cell ends up with no special style applied (matching the rest of the
    # sheet).
    $textRange = $ws.Range("B34:C34")
    $textRange.NumberFormat = "@"

    $ws.Cells.Item(34, 1).Value = 33
    $ws.Cells.Item(34, 2).Value = "2026-02-17"
    $ws.Cells.Item(34, 3).Value = "13:22:43"
    $ws.Cells.Item(34, 4).Value = "MarketMaking"
    $ws.Cells.Item(34, 5).Value = "UP"
    $ws.Cells.Item(34, 6).Value = 0.59
    $ws.Cells.Item(34, 7).Value = 0.7
    $ws.Cells.Item(34, 8).Value = "CLOSED"
    $ws.Cells.Item(34, 9).Value = 18.6441
    $ws.Cells.Item(34, 10).Value = 0.11
    $ws.Cells.Item(34, 11).Value = 98.84999999999999
    $ws.Cells.Item(34, 12).Value = 0
    $ws.Cells.Item(34, 13).Value = 0
    $ws.Cells.Item(34, 14).Value = 0.6
    $ws.Cells.Item(34, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(34, 16).Value = "early_exit"
    $ws.Cells.Item(34, 17).Value = 0.14

    $textRange.ClearFormats()
}

Add-TradeRow34 $wb.Worksheets.Item("All Trades")
Add-TradeRow34 $wb.Worksheets.Item("MarketMaking")
